# GEPEP_calibration.xlsx — calibration update
#
# 1. Update the Weights sheet coefficients (objective is now scaled by the
#    SD of the team variation in the trait).
# 2. Refresh BestBet / Low / High with the 1st iteration of results.
# 3. Remove the stray "20% - Accent6" highlight style from Targets!K2.
# 4. Update selections across sheets and make "High" the active tab.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Weights sheet — new calibrated coefficients
# ---------------------------------------------------------------------
$wsWeights = $wb.Worksheets.Item("Weights")
$wsWeights.Range("B2").Value = 3.0966014953309635
$wsWeights.Range("C2").Value = 1.413120046235661
$wsWeights.Range("D2").Value = 0.0646418641828342
$wsWeights.Range("E2").Value = 0.030508130267892637
$wsWeights.Range("F2").Value = 29.015211341242807
$wsWeights.Range("G2").Value = 7.974570335332701
$wsWeights.Range("H2").Value = 17.457801868409902
$wsWeights.Range("I2").Value = 0.4093757701925073
$wsWeights.Range("J2").Value = 231.42038342001354
$wsWeights.Range("K2").Value = 0.4883772266022255
$wsWeights.Range("J2:L2").Select()

# ---------------------------------------------------------------------
# Targets sheet — drop the highlight style on K2, update selection
# ---------------------------------------------------------------------
$wsTargets = $wb.Worksheets.Item("Targets")
$wsTargets.Range("K2").Style = "Normal"
$wsTargets.Range("L2").Select()

# ---------------------------------------------------------------------
# BestBet sheet — 1st iteration results
# ---------------------------------------------------------------------
$wsBestBet = $wb.Worksheets.Item("BestBet")
$wsBestBet.Range("B2").Value = 3.83410179
$wsBestBet.Range("C2").Value = 18.589628
$wsBestBet.Range("D2").Value = 41.0079669
$wsBestBet.Range("E2").Value = 56392901.1
$wsBestBet.Range("F2").Value = 1.72297182
$wsBestBet.Range("G2").Value = 8.91064517
$wsBestBet.Range("H2").Value = 1.96553236
$wsBestBet.Range("I2").Value = 58.8780614
$wsBestBet.Range("J2").Value = -0.000162630115
$wsBestBet.Range("K2").Value = 1.24978768
$wsBestBet.Range("A2").Select()

# ---------------------------------------------------------------------
# Low sheet — 1st iteration results
# ---------------------------------------------------------------------
$wsLow = $wb.Worksheets.Item("Low")
$wsLow.Range("B2").Value = 2.83410179
$wsLow.Range("C2").Value = 17.589628
$wsLow.Range("D2").Value = 31.0079669
$wsLow.Range("E2").Value = 41392901.1
$wsLow.Range("F2").Value = -2.2770281800000003
$wsLow.Range("G2").Value = 7.9106451700000004
$wsLow.Range("H2").Value = 0.9655323600000001
$wsLow.Range("I2").Value = 53.8780614
$wsLow.Range("J2").Value = -0.00043660271773972607
$wsLow.Range("K2").Value = 0.9997876800000001
$wsLow.Range("J2:M2").Select()

# ---------------------------------------------------------------------
# High sheet — 1st iteration results; becomes the active tab
# ---------------------------------------------------------------------
$wsHigh = $wb.Worksheets.Item("High")
$wsHigh.Range("B2").Value = 4.83410179
$wsHigh.Range("C2").Value = 19.589628
$wsHigh.Range("D2").Value = 51.0079669
$wsHigh.Range("E2").Value = 71392901.1
$wsHigh.Range("F2").Value = 2.7229718199999997
$wsHigh.Range("G2").Value = 9.91064517
$wsHigh.Range("H2").Value = 5.96553236
$wsHigh.Range("I2").Value = 63.8780614
$wsHigh.Range("J2").Value = -0.000025643813630136977
$wsHigh.Range("K2").Value = 1.49978768
$wsHigh.Activate()
$wsHigh.Range("J2:M2").Select()
